$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Capture the formatted content of the "Meta description" paragraph
#    (2nd paragraph of the document) before it gets removed, so its run
#    structure / bold formatting can be reused for the new paragraph that is
#    later introduced near the end of the document.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaFormattedText = $metaPara.Range.FormattedText

# ---------------------------------------------------------------------------
# 2. Insert a new paragraph right before the very last paragraph of the
#    document (the one that used to contain the "Create a feature image..."
#    image-generation prompt). We do this by splitting right after the
#    second-to-last paragraph, resetting the new paragraph to the "Normal"
#    style (so it does not inherit list/heading formatting), and then
#    pasting in the captured formatted text from the "Meta description"
#    paragraph (this preserves the bold run + empty leading run structure).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($count - 1)
$secondToLast.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($count)
$newPara.Range.Style = "Normal"
$newPara.Range.FormattedText = $metaFormattedText

# ---------------------------------------------------------------------------
# 3. Remove the original "Meta description" paragraph entirely.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(2).Range.Delete()

# ---------------------------------------------------------------------------
# 4. The pasted paragraph currently reads:
#        "Meta description: Read our review of the Book of Wealth slot
#         game. Try it for free, featuring Egyptian theme, expanding
#         symbols and free spins bonus."
#    with "Meta description" bold and the remainder (starting at the
#    colon) plain. Trim away everything from the colon onward, leaving
#    just the bold "Meta description" run, then rename that run's text to
#    "Play Book of Wealth Online Slot for Free".
# ---------------------------------------------------------------------------
$count2 = $d.Paragraphs.Count
$newHeadingPara = $d.Paragraphs.Item($count2 - 1)
$headingRange = $newHeadingPara.Range
$trimStart = $headingRange.Start + 16   # length of "Meta description"
$trimRange = $d.Range($trimStart, $headingRange.End)
$trimRange.Delete()

$newHeadingPara2 = $d.Paragraphs.Item($count2 - 1)
$newHeadingPara2.Range.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "Play Book of Wealth Online Slot for Free", 2)

# ---------------------------------------------------------------------------
# 5. Replace the text of the final paragraph (still the italic
#    image-generation prompt) with the new review blurb, keeping its
#    existing (italic) run formatting intact.
# ---------------------------------------------------------------------------
$oldPrompt = "Create a feature image for Book of Wealth Design a fun and exciting feature image for the online slot game, Book of Wealth. The image should be colorful and in a cartoon style. The image should feature a happy Maya warrior with glasses. The warrior should be holding the Book of Wealth and standing in front of a pyramid. The background of the image should have an Egyptian landscape with sand and palm trees. Make sure to include the game logo in the image as well. The image should be eye-catching and attention-grabbing to attract potential players to try out the game."
$newBlurb = "Read our review of the Book of Wealth slot game. Try it for free, featuring Egyptian theme, expanding symbols and free spins bonus."
$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newBlurb, 2)
